$d = $word.ActiveDocument
$r = $d.Range(392, 420)
Write-Host "[392,420) text: [$($r.Text)]"
$r.Delete()
$check = $d.Range(388, 440)
Write-Host "after delete [388,440): [$($check.Text)]"
